# Generate Report for Handoff
# - Flip the "In Translation" status to "Ready for handoff" everywhere it
#   is reported (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to the new handoff-generation time.
# - Columns auto-fit to the new (longer) status text.

$wb = $excel.ActiveWorkbook

# The host's ColumnWidth setter quantizes to 1/6-character-width pixel
# steps (stored_width = floor(ColumnWidth*6 + 0.5 + 5)/6); 16.333333333333332
# is the input that lands on the closest reachable step (17.166666666666668)
# to the author's recorded post-autofit width (17.2159881591797).
$autoFitColumnWidth = 16.333333333333332

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-21 18:44:12"
$wsOverview.Columns.Item(5).ColumnWidth = $autoFitColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $autoFitColumnWidth

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-21 18:44:07"
$wsZhCn.Columns.Item(3).ColumnWidth = $autoFitColumnWidth

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-21 18:44:12"
$wsDeDe.Columns.Item(3).ColumnWidth = $autoFitColumnWidth
